$wb = $excel.ActiveWorkbook

# --- Sheet "Placas" (sheet1): just move the selection to D8 ---
$wsPlacas = $wb.Worksheets.Item("Placas")
$wsPlacas.Activate()
$wsPlacas.Range("D8").Select()

# --- Sheet "Contratos" (sheet2): populate new "Alimentador" values in column C ---
$wsContratos = $wb.Worksheets.Item("Contratos")
$wsContratos.Activate()

$wsContratos.Range("C2").Value = "PCV44482-46KV"
$wsContratos.Range("C3").Value = "PCV32B-4.1KV"
$wsContratos.Range("C4").Value = "PCV44482-46KV"
$wsContratos.Range("C5").Value = "PCV32B-4.1KV"
$wsContratos.Range("C6").Value = "PCV44482-46KV"
$wsContratos.Range("C7").Value = "PCV32B-4.1KV"
$wsContratos.Range("C8").Value = "PCV44482-46KV"
$wsContratos.Range("C9").Value = "PCV32B-4.1KV"

# Adjust column widths as captured in the target file
$wsContratos.Columns.Item(1).ColumnWidth = 11.6
$wsContratos.Columns.Item(3).ColumnWidth = 15.6

# Move the active selection cursor to H7, as recorded in the target file
$wsContratos.Range("H7").Select()
